# lineEstimateTestData.xlsx update
# - Shorten the sample placeholder ID values (drop the "765" suffix) used
#   in the workDetails / adminSanctionDetails / technicalSanctionDetails
#   sheets: AEN765 -> AEN, WIN765 -> WIN, ADN765 -> ADN, TSN765 -> TSN.
# - Move the selection on the workDetails sheet from G2 to E2.

$wb = $excel.ActiveWorkbook

# workDetails: abstractEstimateNumber (E2) and workIdentificationNumber (G2)
$wsWork = $wb.Worksheets.Item("workDetails")
$wsWork.Range("E2").Value = "AEN"
$wsWork.Range("G2").Value = "WIN"
$wsWork.Activate()
[void]$wsWork.Range("E2").Select()

# adminSanctionDetails: administrativeSanctionNumber (B2)
$wsAdmin = $wb.Worksheets.Item("adminSanctionDetails")
$wsAdmin.Range("B2").Value = "ADN"

# technicalSanctionDetails: technicalSanctionNumber (B2)
$wsTech = $wb.Worksheets.Item("technicalSanctionDetails")
$wsTech.Range("B2").Value = "TSN"

# restore the workbook's original active sheet (technicalSanctionDetails)
$wsTech.Activate()
